$wb = $excel.ActiveWorkbook

# This script applies updated market-price / profit figures to the
# "Pandaemonium_Profits" crafting leve tables, one per Disciple of the
# Hand job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW). Values originate
# from a scheduled market-data refresh; no formulas are involved, these
# are plain numeric writes (matching the source workbook, which stores
# only literal <v> values in H:N).

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 277.15
$ws.Range("I33").Value = 232.29411
$ws.Range("J33").Value = 531.3333
$ws.Range("K33").Value = 232.29411
$ws.Range("L33").Value = 531.3333
$ws.Range("M33").Value = -3.294109999999989
$ws.Range("N33").Value = -989.3333
$ws.Range("H51").Value = 3260.2
$ws.Range("I51").Value = 4233.6665
$ws.Range("J51").Value = 1800
$ws.Range("K51").Value = 4233.6665
$ws.Range("L51").Value = 1800
$ws.Range("M51").Value = -3749.6665
$ws.Range("N51").Value = -2768
$ws.Range("H58").Value = 1005.38464
$ws.Range("I58").Value = 181.42857
$ws.Range("J58").Value = 1966.6666
$ws.Range("K58").Value = 544.28571
$ws.Range("L58").Value = 5899.9998
$ws.Range("M58").Value = -394.28571
$ws.Range("N58").Value = -6199.9998
$ws.Range("H100").Value = 1130.7778
$ws.Range("I100").Value = 1089.5714
$ws.Range("J100").Value = 1275
$ws.Range("K100").Value = 1089.5714
$ws.Range("L100").Value = 1275
$ws.Range("M100").Value = -548.5714
$ws.Range("N100").Value = -2357
$ws.Range("H129").Value = 1071.0698
$ws.Range("I129").Value = 325
$ws.Range("J129").Value = 1147.5897
$ws.Range("K129").Value = 975
$ws.Range("L129").Value = 3442.7691
$ws.Range("M129").Value = 4025
$ws.Range("N129").Value = -13442.7691
$ws.Range("H132").Value = 1147.5735
$ws.Range("I132").Value = 1147.5735
$ws.Range("K132").Value = 3442.7205
$ws.Range("M132").Value = -912.7204999999999
$ws.Range("H137").Value = 753334.6
$ws.Range("I137").Value = 2341
$ws.Range("K137").Value = 7023
$ws.Range("M137").Value = -4473
$ws.Range("H138").Value = 5563.984
$ws.Range("I138").Value = 1607.92
$ws.Range("J138").Value = 8166.6577
$ws.Range("K138").Value = 4823.76
$ws.Range("L138").Value = 24499.9731
$ws.Range("M138").Value = 316.2399999999998
$ws.Range("N138").Value = -34779.9731

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18987.045
$ws.Range("I32").Value = 19953.793
$ws.Range("K32").Value = 19953.793
$ws.Range("M32").Value = -19666.793
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 694
$ws.Range("I20").Value = 701.38464
$ws.Range("J20").Value = 655.6
$ws.Range("K20").Value = 701.38464
$ws.Range("L20").Value = 655.6
$ws.Range("M20").Value = -454.38464
$ws.Range("N20").Value = -1149.6
$ws.Range("H22").Value = 345.7143
$ws.Range("I22").Value = 340
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 340
$ws.Range("L22").Value = 350
$ws.Range("M22").Value = -167
$ws.Range("N22").Value = -696
$ws.Range("H122").Value = 56983.285
$ws.Range("J122").Value = 56983.285
$ws.Range("L122").Value = 56983.285
$ws.Range("N122").Value = -66783.285
$ws.Range("H134").Value = 2402.4075
$ws.Range("I134").Value = 2282.55
$ws.Range("J134").Value = 2744.8572
$ws.Range("K134").Value = 6847.650000000001
$ws.Range("L134").Value = 8234.571599999999
$ws.Range("M134").Value = -4312.650000000001
$ws.Range("N134").Value = -13304.5716

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4772.548
$ws.Range("I31").Value = 6059.136
$ws.Range("J31").Value = 3357.3
$ws.Range("K31").Value = 6059.136
$ws.Range("L31").Value = 3357.3
$ws.Range("M31").Value = -5764.136
$ws.Range("N31").Value = -3947.3
$ws.Range("H34").Value = 4772.548
$ws.Range("I34").Value = 6059.136
$ws.Range("J34").Value = 3357.3
$ws.Range("K34").Value = 6059.136
$ws.Range("L34").Value = 3357.3
$ws.Range("M34").Value = -5857.136
$ws.Range("N34").Value = -3761.3
$ws.Range("H59").Value = 22249.5
$ws.Range("J59").Value = 22249.5
$ws.Range("L59").Value = 22249.5
$ws.Range("N59").Value = -24539.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 55.793102
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 55.793102
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 334.758612
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -560.758612
$ws.Range("H38").Value = 51.785713
$ws.Range("I38").Value = 27.5
$ws.Range("J38").Value = 112.5
$ws.Range("K38").Value = 82.5
$ws.Range("L38").Value = 337.5
$ws.Range("M38").Value = 264.5
$ws.Range("N38").Value = -1031.5
$ws.Range("H132").Value = 1608.8518
$ws.Range("I132").Value = 1508.6875
$ws.Range("J132").Value = 1754.5454
$ws.Range("K132").Value = 13578.1875
$ws.Range("L132").Value = 15790.9086
$ws.Range("M132").Value = -11048.1875
$ws.Range("N132").Value = -20850.9086

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2353.6553
$ws.Range("I132").Value = 2288.8262
$ws.Range("J132").Value = 2602.1667
$ws.Range("K132").Value = 6866.4786
$ws.Range("L132").Value = 7806.500100000001
$ws.Range("M132").Value = -4336.4786
$ws.Range("N132").Value = -12866.5001

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5553.4116
$ws.Range("I7").Value = 4954.4614
$ws.Range("J7").Value = 7500
$ws.Range("K7").Value = 4954.4614
$ws.Range("L7").Value = 7500
$ws.Range("M7").Value = -4842.4614
$ws.Range("N7").Value = -7724
$ws.Range("H100").Value = 4632.3335
$ws.Range("I100").Value = 2459.6
$ws.Range("J100").Value = 6184.2856
$ws.Range("K100").Value = 2459.6
$ws.Range("L100").Value = 6184.2856
$ws.Range("M100").Value = -1918.6
$ws.Range("N100").Value = -7266.2856
$ws.Range("H126").Value = 5553.4116
$ws.Range("I126").Value = 4954.4614
$ws.Range("J126").Value = 7500
$ws.Range("K126").Value = 14863.3842
$ws.Range("L126").Value = 22500
$ws.Range("M126").Value = -12393.3842
$ws.Range("N126").Value = -27440
$ws.Range("H136").Value = 3878.6885
$ws.Range("I136").Value = 2162.973
$ws.Range("J136").Value = 6523.75
$ws.Range("K136").Value = 6488.919
$ws.Range("L136").Value = 19571.25
$ws.Range("M136").Value = -3938.919
$ws.Range("N136").Value = -24671.25
$ws.Range("H140").Value = 55248.54
$ws.Range("J140").Value = 55248.54
$ws.Range("L140").Value = 55248.54
$ws.Range("N140").Value = -65608.54000000001
